# Update the "Förändrad" (Changed) date column (C) for every data row
# on the "Avverkningsanmälningar" sheet from 2023-09-08 to 2023-09-09
# (Excel serial date 45177 -> 45178).

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Avverkningsanmälningar") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$newDate = 45178

for ($r = 2; $r -le 261; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
